$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF", formatted like the existing header cells ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# --- Data rows (row -> I, J) ---
$data = @{
    2  = @(7, 8)
    3  = @(7, 7)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(7, 7)
    7  = @(11, 11)
    8  = @(5, 5)
    9  = @(5, 5)
    10 = @(4, 4)
    11 = @(7, 7)
    12 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
